# Add a new "IncidentPrefix" field (header + value) to the Global settings
# sheet, to the right of the existing Language column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# New header + value in column E.
$ws.Range("E1").Value = "IncidentPrefix"
$ws.Range("E2").Value = "Don"

# E2 becomes the new last column in the row, so it picks up the
# "right border" formatting that D2 (the old last column) used to have.
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# D2 is no longer the last column, so it now gets the plain "interior
# column" formatting (no right border), same as the other inner cells.
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

# Size the new column to fit its contents.
$ws.Columns.Item(5).AutoFit()

# Move the selection down past the two rows of data.
$ws.Range("A3").Select()

# Restore the originally active sheet/tab.
$wsOther = $wb.Worksheets.Item("CreateIncident")
$wsOther.Activate()
